$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, shifting old F (and everything after) one to the right.
$ws.Columns("F:F").Insert()

# Best-effort: carry column F's width/alignment over from column E so the
# new column looks consistent with its neighbour.
$eCol = $ws.Columns("E:E")
$fCol = $ws.Columns("F:F")
$fCol.ColumnWidth = $eCol.ColumnWidth
$fCol.HorizontalAlignment = $eCol.HorizontalAlignment
$fCol.VerticalAlignment = $eCol.VerticalAlignment
$fCol.WrapText = $eCol.WrapText

# Introduce the brand-new strings in the same order the author first typed
# them (this determines their position in xl/sharedStrings.xml).
$ws.Range("A2").Value = "Activiteit"
$ws.Range("A1").Value = "[Activiteiten]"
$ws.Range("C1").Value = "actRole"
$ws.Range("F1").Value = "actRandcdx"
$ws.Range("G1").Value = "actVoorschrift"
$ws.Range("B2").Value = "ActNaam"

# --- Row 1 (header row 1) ---
$ws.Range("B1").Value = "actNaam"
$ws.Range("D1").Value = "actPrecdx"
$ws.Range("E1").Value = "actPostcdx"

# --- Row 2 (header row 2) ---
$ws.Range("C2").Value = "Role"
$ws.Range("D2").Value = "Expressie"
$ws.Range("E2").Value = "Expressie"
$ws.Range("F2").Value = "Expressie"
$ws.Range("G2").Value = "Handelingsvoorschrift"

# --- Row 3 ---
$ws.Range("A3").Value = "Actie_0001"
$ws.Range("B3").Value = "registreer burger"
$ws.Range("C3").Value = "admin"
$ws.Range("D3").Value = "deze expressie is altijd waar"
$ws.Range("E3").Value = "burger is geregisteerd"
$ws.Range("G3").Value = "burger moet worden opgenomen in het register van geautoriseerde voters"

# --- Row 4 ---
$ws.Range("A4").Value = "Actie_0002"
$ws.Range("B4").Value = "marietje stemt"
$ws.Range("C4").Value = "burger"
$ws.Range("D4").Value = "marietje heeft niet gestemd"
$ws.Range("E4").Value = "marietje heeft gestemd OF de stemming is niet (meer) opengesteld"
$ws.Range("G4").Value = "marietje brengt haar stem uit (en aanvullende instructies)"

# --- Row 5 ---
$ws.Range("A5").Value = "Actie_0003"
$ws.Range("B5").Value = "jantje stemt"
$ws.Range("C5").Value = "burger"
$ws.Range("D5").Value = "jantje is geregistreerd EN jantje heeft niet gestemd"
$ws.Range("E5").Value = "jantje heeft gestemd OF de stemming is niet (meer) opengesteld"
$ws.Range("G5").Value = "jantje moet naar het stembureau om zijn stem uit te brengen"

# --- Row 6 ---
$ws.Range("A6").Value = "Actie_0004"
$ws.Range("B6").Value = "klaas is kandidaat"
$ws.Range("C6").Value = "admin"
$ws.Range("D6").Value = "klaas heeft kwalificaties voor het kandidaatschap"
$ws.Range("E6").Value = "klaas is geregistreerd als kandidaat OF de kandidaatstellingstermijn is verlopen"
$ws.Range("G6").Value = "admin moet klaas op de kandidatenlijst zetten"

# The newly-inserted column F has no content on the plain data rows (3-6);
# make sure no stray formatting / placeholder cell lingers there, matching
# a plain column insert with nothing ever typed into those cells.
$ws.Range("F3:F6").Clear()
